$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = -0.009549999999999999
$ws.Range("E2").Value = -0.159825
$ws.Range("G2").Value = 0.263308302044222
$ws.Range("H2").Value = 0.263308302044222
$ws.Range("I2").Value = 0.2381309970796829
$ws.Range("J2").Value = 0.219914595935892
$ws.Range("K2").Value = 78.51000000000001
$ws.Range("L2").Value = 0.1637672090112641
$ws.Range("M2").Value = 73.15299999999999
$ws.Range("N2").Value = 0.04538027295285359
$ws.Range("O2").Value = 0.9317666539294356
$ws.Range("P2").Value = 72.41
$ws.Range("Q2").Value = 0.04491935483870967
$ws.Range("R2").Value = 0.9223028913514201
$ws.Range("S2").Value = 0.7430000000000003
$ws.Range("T2").Value = 0.01015679466324006
$ws.Range("U2").Value = 66.90000000000001
$ws.Range("V2").Value = 0.04150124069478908
$ws.Range("W2").Value = 0.005963855421686747
$ws.Range("X2").Value = 0.0482493014422786
$ws.Range("Y2").Value = -0.04228544602059185
$ws.Range("Z2").Value = 0.6040826612903226
$ws.Range("AA2").Value = 0.03377061956098195
$ws.Range("AB2").Value = 0.04787080373776594
$ws.Range("AC2").Value = -0.01410018417678399
$ws.Range("AD2").Value = 39.01000000000001
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 39.01000000000001
$ws.Range("AG2").Value = -27.89
$ws.Range("AH2").Value = 0.02362796106625642
$ws.Range("AI2").Value = 0.03897453317481093
$ws.Range("AJ2").Value = -0.01760610058644917
$ws.Range("AK2").Value = -0.02986049399899359
$ws.Range("AL2").Value = 1.84
$ws.Range("AM2").Value = 1.84
$ws.Range("AN2").Value = 0.332197905134974
$ws.Range("AO2").Value = 62.04347826086956
$ws.Range("AP2").Value = -0.2375031933918079
$ws.Range("AQ2").Value = 62.04347826086956
$ws.Range("D3").Value = 0.04940000000000001
$ws.Range("E3").Value = -0.00465
$ws.Range("G3").Value = 0.3470708446866485
$ws.Range("H3").Value = 0.3470708446866485
$ws.Range("I3").Value = 0.3480926430517711
$ws.Range("J3").Value = 0.2682081009451635
$ws.Range("K3").Value = 78.90000000000001
$ws.Range("L3").Value = 0.2687329700272479
$ws.Range("M3").Value = 68.09999999999999
$ws.Range("N3").Value = 0.05006984780530843
$ws.Range("O3").Value = 0.8631178707224333
$ws.Range("P3").Value = 68.09999999999999
$ws.Range("Q3").Value = 0.05006984780530843
$ws.Range("R3").Value = 0.8631178707224333
$ws.Range("U3").Value = 24.6
$ws.Range("V3").Value = 0.01808690537460481
$ws.Range("W3").Value = 0.1727611123275674
$ws.Range("X3").Value = 0.04804246279368089
$ws.Range("Y3").Value = 0.1247186495338865
$ws.Range("Z3").Value = 0.7018885967009324
$ws.Range("AA3").Value = 0.1882522075962228
$ws.Range("AB3").Value = 0.04784242501435481
$ws.Range("AC3").Value = 0.140409782581868
$ws.Range("AD3").Value = 10.8
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 10.8
$ws.Range("AG3").Value = -13.8
$ws.Range("AH3").Value = 0.00787803632650084
$ws.Range("AI3").Value = 0.02478769795731008
$ws.Range("AJ3").Value = -0.01025031568001188
$ws.Range("AK3").Value = -0.03356847482364389
$ws.Range("AL3").Value = 0.448
$ws.Range("AM3").Value = 0.448
$ws.Range("AN3").Value = 0.1046511627906977
$ws.Range("AO3").Value = 228.125
$ws.Range("AP3").Value = -0.1337209302325582
$ws.Range("AQ3").Value = 228.125
$ws.Range("D4").Value = -0.009549999999999999
$ws.Range("E4").Value = -0.315
$ws.Range("G4").Value = 0.2149080348499516
$ws.Range("H4").Value = 0.2149080348499516
$ws.Range("I4").Value = 0.1006776379477251
$ws.Range("J4").Value = 0.1006776379477251
$ws.Range("K4").Value = 1.98
$ws.Range("L4").Value = 0.0191674733785092
$ws.Range("U4").Value = 25.2
$ws.Range("V4").Value = 0.1299638989169675
$ws.Range("W4").Value = 0.005963855421686747
$ws.Range("X4").Value = 0.0482493014422786
$ws.Range("Y4").Value = -0.04228544602059185
$ws.Range("Z4").Value = 0.3354331731393688
$ws.Range("AA4").Value = 0.03377061956098195
$ws.Range("AB4").Value = 0.04787080373776594
$ws.Range("AC4").Value = -0.01410018417678399
$ws.Range("AD4").Value = 2.91
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 2.91
$ws.Range("AG4").Value = -22.29
$ws.Range("AH4").Value = 0.01478583405314771
$ws.Range("AI4").Value = 0.006472275972509508
$ws.Range("AJ4").Value = -0.1298875356913932
$ws.Range("AK4").Value = -0.05251996889799958
$ws.Range("AL4").Value = 0.152
$ws.Range("AM4").Value = 0.152
$ws.Range("AN4").Value = 0.2404958677685951
$ws.Range("AO4").Value = 68.42105263157895
$ws.Range("AP4").Value = -1.842148760330578
$ws.Range("AQ4").Value = 68.42105263157895
$ws.Range("D5").Value = -0.0113
$ws.Range("G5").Value = 0.02581818181818182
$ws.Range("H5").Value = 0.02581818181818182
$ws.Range("I5").Value = 0.01890909090909091
$ws.Range("J5").Value = 0.01890909090909091
$ws.Range("K5").Value = -2.37
$ws.Range("L5").Value = -0.02872727272727273
$ws.Range("M5").Value = 5.053
$ws.Range("N5").Value = 0.08712068965517242
$ws.Range("O5").Value = -2.132067510548523
$ws.Range("P5").Value = 4.31
$ws.Range("Q5").Value = 0.07431034482758619
$ws.Range("R5").Value = -1.818565400843882
$ws.Range("S5").Value = 0.7430000000000003
$ws.Range("T5").Value = 0.1470413615673858
$ws.Range("U5").Value = 17.1
$ws.Range("V5").Value = 0.2948275862068966
$ws.Range("W5").Value = -0.03424855491329479
$ws.Range("X5").Value = 0.060576809311195
$ws.Range("Y5").Value = -0.09482536422448978
$ws.Range("Z5").Value = 1.225126225126225
$ws.Range("AA5").Value = 0.02316602316602316
$ws.Range("AB5").Value = 0.04905781374418811
$ws.Range("AC5").Value = -0.02589179057816494
$ws.Range("AD5").Value = 25.3
$ws.Range("AF5").Value = 25.3
$ws.Range("AG5").Value = 8.199999999999999
$ws.Range("AH5").Value = 0.3037214885954382
$ws.Range("AI5").Value = 0.2188581314878893
$ws.Range("AJ5").Value = 0.1238670694864048
$ws.Range("AK5").Value = 0.083248730964467
$ws.Range("AL5").Value = 1.24
$ws.Range("AM5").Value = 1.24
$ws.Range("AN5").Value = 11.87793427230047
$ws.Range("AO5").Value = 1.258064516129032
$ws.Range("AP5").Value = 3.849765258215962
$ws.Range("AQ5").Value = 1.258064516129032
